$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Productos")

# Update the three product rows with new product names / quantities
$ws.Range("B2").Value = "Encuentra Tu Persona Vitamina PLANETA 3028371"

$ws.Range("B3").Value = "Primeros Mil Días Del Bebe GRIJALBO 1303320"
$ws.Range("C3").Value = "4.00 un"

$ws.Range("B4").Value = "Cree En Ti DIANA 3026463"
